$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Is Active" (F2) was a text string "False"; push up a real boolean value instead.
$ws.Range("F2").Value = $False
